$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 500116.7
$ws.Range("I9").Value = 114.625
$ws.Range("K9").Value = 114.625
$ws.Range("M9").Value = 54.375
$ws.Range("H12").Value = 93.2
$ws.Range("I12").Value = 81.44444
$ws.Range("J12").Value = 199
$ws.Range("K12").Value = 81.44444
$ws.Range("L12").Value = 199
$ws.Range("M12").Value = 88.55556
$ws.Range("N12").Value = -539
$ws.Range("H15").Value = 1711.1621
$ws.Range("I15").Value = 1711.1621
$ws.Range("K15").Value = 5133.4863
$ws.Range("M15").Value = -4964.4863
$ws.Range("H18").Value = 548.2
$ws.Range("I18").Value = 548.2
$ws.Range("K18").Value = 548.2
$ws.Range("M18").Value = -264.2
$ws.Range("H28").Value = 195.7
$ws.Range("I28").Value = 195.7
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 195.7
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 289.3
$ws.Range("N28").Value = ""
$ws.Range("H33").Value = 300.78262
$ws.Range("I33").Value = 321.66666
$ws.Range("K33").Value = 321.66666
$ws.Range("M33").Value = -92.66665999999998
$ws.Range("H40").Value = 7111.875
$ws.Range("I40").Value = 8683
$ws.Range("J40").Value = 2398.5
$ws.Range("K40").Value = 8683
$ws.Range("L40").Value = 2398.5
$ws.Range("M40").Value = -8508
$ws.Range("N40").Value = -2748.5
$ws.Range("H43").Value = 2017.8
$ws.Range("I43").Value = 1539.25
$ws.Range("J43").Value = 2191.818
$ws.Range("K43").Value = 1539.25
$ws.Range("L43").Value = 2191.818
$ws.Range("M43").Value = -1470.25
$ws.Range("N43").Value = -2329.818
$ws.Range("H62").Value = 5977
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H64").Value = 6405.769
$ws.Range("I64").Value = 5030.5557
$ws.Range("J64").Value = 9500
$ws.Range("K64").Value = 5030.5557
$ws.Range("L64").Value = 9500
$ws.Range("M64").Value = -4782.5557
$ws.Range("N64").Value = -9996
$ws.Range("H65").Value = 5977
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H67").Value = 6405.769
$ws.Range("I67").Value = 5030.5557
$ws.Range("J67").Value = 9500
$ws.Range("K67").Value = 5030.5557
$ws.Range("L67").Value = 9500
$ws.Range("M67").Value = -4172.5557
$ws.Range("N67").Value = -11216
$ws.Range("H70").Value = 126499.875
$ws.Range("I70").Value = 200320
$ws.Range("K70").Value = 600960
$ws.Range("M70").Value = -600690
$ws.Range("H73").Value = 126499.875
$ws.Range("I73").Value = 200320
$ws.Range("K73").Value = 600960
$ws.Range("M73").Value = -600024
$ws.Range("H74").Value = 9326.666999999999
$ws.Range("I74").Value = 8995
$ws.Range("K74").Value = 8995
$ws.Range("M74").Value = -8059
$ws.Range("H76").Value = 11986.4
$ws.Range("I76").Value = 11984.889
$ws.Range("K76").Value = 11984.889
$ws.Range("M76").Value = -11669.889
$ws.Range("H77").Value = 9326.666999999999
$ws.Range("I77").Value = 8995
$ws.Range("K77").Value = 44975
$ws.Range("M77").Value = -40295
$ws.Range("H79").Value = 11986.4
$ws.Range("I79").Value = 11984.889
$ws.Range("K79").Value = 11984.889
$ws.Range("M79").Value = -10892.889
$ws.Range("H98").Value = 678.7826
$ws.Range("I98").Value = 678.7826
$ws.Range("K98").Value = 678.7826
$ws.Range("M98").Value = 819.2174
$ws.Range("H106").Value = 5348.125
$ws.Range("I106").Value = 5348.125
$ws.Range("K106").Value = 5348.125
$ws.Range("M106").Value = -4717.125
$ws.Range("H107").Value = 613.6429000000001
$ws.Range("I107").Value = 601
$ws.Range("J107").Value = 689.5
$ws.Range("K107").Value = 601
$ws.Range("L107").Value = 689.5
$ws.Range("M107").Value = 1319
$ws.Range("N107").Value = -4529.5
$ws.Range("H111").Value = 3999.5
$ws.Range("I111").Value = 3999.5
$ws.Range("K111").Value = 11998.5
$ws.Range("M111").Value = -8931.5
$ws.Range("H112").Value = 2010.3529
$ws.Range("J112").Value = 2241.3076
$ws.Range("L112").Value = 6723.9228
$ws.Range("N112").Value = -8939.9228
$ws.Range("H113").Value = 12249.75
$ws.Range("I113").Value = 12249.75
$ws.Range("K113").Value = 12249.75
$ws.Range("M113").Value = -8995.75
$ws.Range("H122").Value = 678.7826
$ws.Range("I122").Value = 678.7826
$ws.Range("K122").Value = 2036.3478
$ws.Range("M122").Value = 413.6522
$ws.Range("H125").Value = 1785.1428
$ws.Range("I125").Value = 1624.5
$ws.Range("K125").Value = 14620.5
$ws.Range("M125").Value = -12160.5
$ws.Range("H137").Value = 3819.85
$ws.Range("I137").Value = 3481.2
$ws.Range("J137").Value = 3932.7334
$ws.Range("K137").Value = 10443.6
$ws.Range("L137").Value = 11798.2002
$ws.Range("M137").Value = -7893.599999999999
$ws.Range("N137").Value = -16898.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 13050
$ws.Range("I3").Value = 25500
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 25500
$ws.Range("L3").Value = 600
$ws.Range("M3").Value = -25385
$ws.Range("N3").Value = -830
$ws.Range("H4").Value = 365.5
$ws.Range("I4").Value = 359.75
$ws.Range("K4").Value = 359.75
$ws.Range("M4").Value = -243.75
$ws.Range("H32").Value = 21406.51
$ws.Range("I32").Value = 8336.727999999999
$ws.Range("J32").Value = 52213.855
$ws.Range("K32").Value = 8336.727999999999
$ws.Range("L32").Value = 52213.855
$ws.Range("M32").Value = -8049.727999999999
$ws.Range("N32").Value = -52787.855
$ws.Range("H45").Value = 1373.6
$ws.Range("I45").Value = 1469.5
$ws.Range("J45").Value = 990
$ws.Range("K45").Value = 1469.5
$ws.Range("L45").Value = 990
$ws.Range("M45").Value = -1092.5
$ws.Range("N45").Value = -1744
$ws.Range("H61").Value = 3269.8235
$ws.Range("I61").Value = 3141.5833
$ws.Range("J61").Value = 3577.6
$ws.Range("K61").Value = 3141.5833
$ws.Range("L61").Value = 3577.6
$ws.Range("M61").Value = -2929.5833
$ws.Range("N61").Value = -4001.6
$ws.Range("H63").Value = 2561.4905
$ws.Range("I63").Value = 2545.5334
$ws.Range("K63").Value = 2545.5334
$ws.Range("M63").Value = -1859.5334
$ws.Range("H66").Value = 2561.4905
$ws.Range("I66").Value = 2545.5334
$ws.Range("K66").Value = 12727.667
$ws.Range("M66").Value = -9295.666999999999
$ws.Range("H102").Value = 3583.5833
$ws.Range("I102").Value = 3300.4
$ws.Range("J102").Value = 4999.5
$ws.Range("K102").Value = 3300.4
$ws.Range("L102").Value = 4999.5
$ws.Range("M102").Value = -1678.4
$ws.Range("N102").Value = -8243.5
$ws.Range("H122").Value = 2195.923
$ws.Range("I122").Value = 2042.0294
$ws.Range("K122").Value = 6126.0882
$ws.Range("M122").Value = -3676.0882
$ws.Range("H132").Value = 2453.182
$ws.Range("I132").Value = 2001.5
$ws.Range("K132").Value = 6004.5
$ws.Range("M132").Value = -3474.5
$ws.Range("H135").Value = 79299.60000000001
$ws.Range("J135").Value = 79299.60000000001
$ws.Range("L135").Value = 79299.60000000001
$ws.Range("N135").Value = -89439.60000000001
$ws.Range("H136").Value = 3269.8235
$ws.Range("I136").Value = 3141.5833
$ws.Range("J136").Value = 3577.6
$ws.Range("K136").Value = 9424.749899999999
$ws.Range("L136").Value = 10732.8
$ws.Range("M136").Value = -6874.749899999999
$ws.Range("N136").Value = -15832.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3179.4666
$ws.Range("I86").Value = 3179.4666
$ws.Range("K86").Value = 3179.4666
$ws.Range("M86").Value = -2056.4666
$ws.Range("H89").Value = 3179.4666
$ws.Range("I89").Value = 3179.4666
$ws.Range("K89").Value = 15897.333
$ws.Range("M89").Value = -10281.333
$ws.Range("H105").Value = 5001.05
$ws.Range("I105").Value = 4942.5557
$ws.Range("J105").Value = 5527.5
$ws.Range("K105").Value = 4942.5557
$ws.Range("L105").Value = 5527.5
$ws.Range("M105").Value = -3195.5557
$ws.Range("N105").Value = -9021.5
$ws.Range("H134").Value = 5621.2383
$ws.Range("I134").Value = 3543.5386
$ws.Range("K134").Value = 10630.6158
$ws.Range("M134").Value = -8095.6158
$ws.Range("H135").Value = 58750
$ws.Range("J135").Value = 58750
$ws.Range("L135").Value = 58750
$ws.Range("N135").Value = -68890
$ws.Range("H137").Value = 58500
$ws.Range("J137").Value = 58500
$ws.Range("L137").Value = 58500
$ws.Range("N137").Value = -68700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 111287.11
$ws.Range("J7").Value = 297.5
$ws.Range("L7").Value = 297.5
$ws.Range("N7").Value = -523.5
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = ""
$ws.Range("H31").Value = 5799.8726
$ws.Range("I31").Value = 4404.476
$ws.Range("J31").Value = 6926.923
$ws.Range("K31").Value = 4404.476
$ws.Range("L31").Value = 6926.923
$ws.Range("M31").Value = -4109.476
$ws.Range("N31").Value = -7516.923
$ws.Range("H34").Value = 5799.8726
$ws.Range("I34").Value = 4404.476
$ws.Range("J34").Value = 6926.923
$ws.Range("K34").Value = 4404.476
$ws.Range("L34").Value = 6926.923
$ws.Range("M34").Value = -4202.476
$ws.Range("N34").Value = -7330.923
$ws.Range("H55").Value = 25357.334
$ws.Range("I55").Value = 73
$ws.Range("J55").Value = 37999.5
$ws.Range("K55").Value = 73
$ws.Range("L55").Value = 37999.5
$ws.Range("M55").Value = 242
$ws.Range("N55").Value = -38629.5
$ws.Range("H58").Value = 6098.778
$ws.Range("I58").Value = 6098.778
$ws.Range("K58").Value = 6098.778
$ws.Range("M58").Value = -5895.778
$ws.Range("H62").Value = 5210.75
$ws.Range("I62").Value = 5240.857
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 5240.857
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4616.857
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 5210.75
$ws.Range("I65").Value = 5240.857
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 26204.285
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -23084.285
$ws.Range("N65").Value = -31240
$ws.Range("H88").Value = 15135.286
$ws.Range("J88").Value = 15135.286
$ws.Range("L88").Value = 15135.286
$ws.Range("N88").Value = -15947.286
$ws.Range("H91").Value = 15135.286
$ws.Range("J91").Value = 15135.286
$ws.Range("L91").Value = 15135.286
$ws.Range("N91").Value = -17943.286
$ws.Range("H99").Value = 2834.4
$ws.Range("I99").Value = 2834.4
$ws.Range("K99").Value = 2834.4
$ws.Range("M99").Value = -1336.4
$ws.Range("H105").Value = 4764
$ws.Range("I105").Value = 4629.5557
$ws.Range("J105").Value = 4965.6665
$ws.Range("K105").Value = 4629.5557
$ws.Range("L105").Value = 4965.6665
$ws.Range("M105").Value = -2882.5557
$ws.Range("N105").Value = -8459.666499999999
$ws.Range("H126").Value = 2834.4
$ws.Range("I126").Value = 2834.4
$ws.Range("K126").Value = 8503.200000000001
$ws.Range("M126").Value = -6033.200000000001
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = ""
$ws.Range("H134").Value = 1166.9546
$ws.Range("I134").Value = 1188.1666
$ws.Range("K134").Value = 3564.4998
$ws.Range("M134").Value = -1029.4998
$ws.Range("H136").Value = 6098.778
$ws.Range("I136").Value = 6098.778
$ws.Range("K136").Value = 18296.334
$ws.Range("M136").Value = -15746.334
$ws.Range("H137").Value = 95000
$ws.Range("J137").Value = 95000
$ws.Range("L137").Value = 95000
$ws.Range("N137").Value = -105200
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""
$ws.Range("H140").Value = 99000
$ws.Range("J140").Value = 99000
$ws.Range("L140").Value = 99000
$ws.Range("N140").Value = -109360
$ws.Range("H141").Value = 469790.4
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 469790.4
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 469790.4
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -480150.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 144.06667
$ws.Range("I2").Value = 122.71429
$ws.Range("J2").Value = 162.75
$ws.Range("K2").Value = 736.28574
$ws.Range("L2").Value = 976.5
$ws.Range("M2").Value = -623.28574
$ws.Range("N2").Value = -1202.5
$ws.Range("H14").Value = 274.5
$ws.Range("I14").Value = 274.5
$ws.Range("K14").Value = 823.5
$ws.Range("M14").Value = -650.5
$ws.Range("H17").Value = 2698.5
$ws.Range("J17").Value = 4900
$ws.Range("L17").Value = 14700
$ws.Range("N17").Value = -15038
$ws.Range("H75").Value = 2750.4285
$ws.Range("I75").Value = 1830.75
$ws.Range("J75").Value = 3118.3
$ws.Range("K75").Value = 5492.25
$ws.Range("L75").Value = 9354.900000000001
$ws.Range("M75").Value = -4494.25
$ws.Range("N75").Value = -11350.9
$ws.Range("H78").Value = 2750.4285
$ws.Range("I78").Value = 1830.75
$ws.Range("J78").Value = 3118.3
$ws.Range("K78").Value = 16476.75
$ws.Range("L78").Value = 28064.7
$ws.Range("M78").Value = -11484.75
$ws.Range("N78").Value = -38048.7
$ws.Range("H102").Value = 2580
$ws.Range("I102").Value = 2580
$ws.Range("K102").Value = 7740
$ws.Range("M102").Value = -5306
$ws.Range("H113").Value = 425.9375
$ws.Range("I113").Value = 286
$ws.Range("J113").Value = 435.26666
$ws.Range("K113").Value = 858
$ws.Range("L113").Value = 1305.79998
$ws.Range("M113").Value = 1312
$ws.Range("N113").Value = -5645.79998
$ws.Range("H121").Value = 68805.2
$ws.Range("I121").Value = 111779.89
$ws.Range("J121").Value = 4343.1665
$ws.Range("K121").Value = 335339.67
$ws.Range("L121").Value = 13029.4995
$ws.Range("M121").Value = -334029.67
$ws.Range("N121").Value = -15649.4995
$ws.Range("H132").Value = 1057.091
$ws.Range("I132").Value = 1057.091
$ws.Range("K132").Value = 9513.819
$ws.Range("M132").Value = -6983.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 242.5625
$ws.Range("I2").Value = 242.5625
$ws.Range("K2").Value = 242.5625
$ws.Range("M2").Value = -129.5625
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""
$ws.Range("H23").Value = 29995
$ws.Range("J23").Value = 29995
$ws.Range("L23").Value = 29995
$ws.Range("N23").Value = -30441
$ws.Range("H24").Value = 30000
$ws.Range("I24").Value = 30000
$ws.Range("K24").Value = 30000
$ws.Range("M24").Value = -29827
$ws.Range("H80").Value = 17813.947
$ws.Range("I80").Value = 14496.667
$ws.Range("J80").Value = 18435.938
$ws.Range("K80").Value = 14496.667
$ws.Range("L80").Value = 18435.938
$ws.Range("M80").Value = -13498.667
$ws.Range("N80").Value = -20431.938
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H83").Value = 17813.947
$ws.Range("I83").Value = 14496.667
$ws.Range("J83").Value = 18435.938
$ws.Range("K83").Value = 72483.33499999999
$ws.Range("L83").Value = 92179.68999999999
$ws.Range("M83").Value = -67491.33499999999
$ws.Range("N83").Value = -102163.69
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H97").Value = 1174.5
$ws.Range("I97").Value = 1174.5
$ws.Range("K97").Value = 1174.5
$ws.Range("M97").Value = -678.5
$ws.Range("H102").Value = 28530.854
$ws.Range("I102").Value = 30372.773
$ws.Range("J102").Value = 9497.666999999999
$ws.Range("K102").Value = 30372.773
$ws.Range("L102").Value = 9497.666999999999
$ws.Range("M102").Value = -28750.773
$ws.Range("N102").Value = -12741.667
$ws.Range("H107").Value = 780.8889
$ws.Range("I107").Value = 505
$ws.Range("J107").Value = 1746.5
$ws.Range("K107").Value = 505
$ws.Range("L107").Value = 1746.5
$ws.Range("M107").Value = 1415
$ws.Range("N107").Value = -5586.5
$ws.Range("H113").Value = 2047.45
$ws.Range("I113").Value = 1997.3158
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1997.3158
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 172.6841999999999
$ws.Range("N113").Value = -7340
$ws.Range("H122").Value = 57671.89
$ws.Range("I122").Value = 144570.86
$ws.Range("K122").Value = 433712.58
$ws.Range("M122").Value = -431262.58
$ws.Range("H126").Value = 1498.5
$ws.Range("I126").Value = 998
$ws.Range("K126").Value = 2994
$ws.Range("M126").Value = -524
$ws.Range("H132").Value = 2580.5862
$ws.Range("I132").Value = 2646.125
$ws.Range("K132").Value = 7938.375
$ws.Range("M132").Value = -5408.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8464.75
$ws.Range("I7").Value = 7168.7
$ws.Range("K7").Value = 7168.7
$ws.Range("M7").Value = -7056.7
$ws.Range("H22").Value = 1121.909
$ws.Range("I22").Value = 996.5
$ws.Range("J22").Value = 1149.7778
$ws.Range("K22").Value = 996.5
$ws.Range("L22").Value = 1149.7778
$ws.Range("M22").Value = -701.5
$ws.Range("N22").Value = -1739.7778
$ws.Range("H27").Value = 1121.909
$ws.Range("I27").Value = 996.5
$ws.Range("J27").Value = 1149.7778
$ws.Range("K27").Value = 996.5
$ws.Range("L27").Value = 1149.7778
$ws.Range("M27").Value = -889.5
$ws.Range("N27").Value = -1363.7778
$ws.Range("H40").Value = 3061.6072
$ws.Range("I40").Value = 2514.3809
$ws.Range("J40").Value = 4703.2856
$ws.Range("K40").Value = 2514.3809
$ws.Range("L40").Value = 4703.2856
$ws.Range("M40").Value = -2378.3809
$ws.Range("N40").Value = -4975.2856
$ws.Range("H55").Value = 550.125
$ws.Range("I55").Value = 618.5
$ws.Range("J55").Value = 345
$ws.Range("K55").Value = 618.5
$ws.Range("L55").Value = 345
$ws.Range("M55").Value = -445.5
$ws.Range("N55").Value = -691
$ws.Range("H62").Value = 34990
$ws.Range("J62").Value = 34990
$ws.Range("L62").Value = 34990
$ws.Range("N62").Value = -36238
$ws.Range("H65").Value = 34990
$ws.Range("J65").Value = 34990
$ws.Range("L65").Value = 104970
$ws.Range("N65").Value = -111210
$ws.Range("H93").Value = 2524.65
$ws.Range("I93").Value = 2453.375
$ws.Range("K93").Value = 2453.375
$ws.Range("M93").Value = -1205.375
$ws.Range("H122").Value = 2400
$ws.Range("I122").Value = 2400
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7200
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4750
$ws.Range("N122").Value = ""
$ws.Range("H126").Value = 8464.75
$ws.Range("I126").Value = 7168.7
$ws.Range("K126").Value = 21506.1
$ws.Range("M126").Value = -19036.1
$ws.Range("H132").Value = 2313.4722
$ws.Range("I132").Value = 2536.6
$ws.Range("J132").Value = 2034.5625
$ws.Range("K132").Value = 7609.799999999999
$ws.Range("L132").Value = 6103.6875
$ws.Range("M132").Value = -5079.799999999999
$ws.Range("N132").Value = -11163.6875
$ws.Range("H136").Value = 3145.8
$ws.Range("I136").Value = 2398.818
$ws.Range("J136").Value = 5200
$ws.Range("K136").Value = 7196.454000000001
$ws.Range("L136").Value = 15600
$ws.Range("M136").Value = -4646.454000000001
$ws.Range("N136").Value = -20700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 18000
$ws.Range("I8").Value = 18000
$ws.Range("K8").Value = 18000
$ws.Range("M8").Value = -17860
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = ""
$ws.Range("H100").Value = 2934.3333
$ws.Range("I100").Value = 2934.3333
$ws.Range("K100").Value = 5868.6666
$ws.Range("M100").Value = -5327.6666
$ws.Range("H105").Value = 22871.334
$ws.Range("J105").Value = 22871.334
$ws.Range("L105").Value = 22871.334
$ws.Range("N105").Value = -29859.334
$ws.Range("H107").Value = 825.25
$ws.Range("I107").Value = 414.6
$ws.Range("J107").Value = 1509.6666
$ws.Range("K107").Value = 1243.8
$ws.Range("L107").Value = 4528.9998
$ws.Range("M107").Value = 676.1999999999998
$ws.Range("N107").Value = -8368.9998
$ws.Range("H125").Value = 63629
$ws.Range("J125").Value = 63629
$ws.Range("L125").Value = 63629
$ws.Range("N125").Value = -73469
$ws.Range("H130").Value = 58999
$ws.Range("J130").Value = 58999
$ws.Range("L130").Value = 58999
$ws.Range("N130").Value = -69039
$ws.Range("H132").Value = 3729.5
$ws.Range("I132").Value = 3699.5557
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 11098.6671
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -8568.667099999999
$ws.Range("N132").Value = -17057
$ws.Range("H136").Value = 2105.353
$ws.Range("I136").Value = 2399.3
$ws.Range("J136").Value = 1685.4286
$ws.Range("K136").Value = 7197.900000000001
$ws.Range("L136").Value = 5056.2858
$ws.Range("M136").Value = -4647.900000000001
$ws.Range("N136").Value = -10156.2858
